$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SCORE (D2) and CLUSTER MUST HAVE SCORE (G2) are stored as text in the
# source workbook (e.g. '26.80' must keep the trailing zero), so force the
# cells to text format before writing the new values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("D2").Value = '26.80'
$ws.Range("G2").Value = '27.77'

# CLUSTER MUST HAVE MATCH column (F2:F6) - same keywords, new order
$ws.Range("F2").Value = 'project : 1'
$ws.Range("F3").Value = 'analysis : 1'
$ws.Range("F4").Value = 'git : 1'
$ws.Range("F5").Value = 'sql : 1'
$ws.Range("F6").Value = 'databases : 2'

# CLUSTER GOOD TO HAVE MATCH column (H2:H3) - same keywords, new order
$ws.Range("H2").Value = 'regulations : 2'
$ws.Range("H3").Value = 'process : 1'

# CV KEYWORDS column (E2:E218) - refreshed keyword analysis, now
# including 2-word skills (e.g. 'oracle dba', 'system administrator')
$keywords = @(
  'image processing : 1',
  'communications : 1',
  'electronics : 1',
  'contract : 23',
  'energy : 2',
  'plc : 4',
  'tv : 2',
  'financial services : 3',
  'international : 2',
  'network : 8',
  'systems : 12',
  'telecommunications : 3',
  'oracle dba : 53',
  'dba : 73',
  'consultant : 8',
  'unix : 7',
  'system administrator : 7',
  'administrator : 5',
  'analyst : 3',
  'manager : 12',
  'acquisition : 1',
  'admin : 1',
  'oracle : 130',
  'oracle applications : 3',
  'middleware : 4',
  'oem : 19',
  'monitor : 1',
  'oracle databases : 8',
  'databases : 28',
  'patching : 8',
  'tuning : 30',
  'monitoring : 16',
  'rman : 20',
  'testing : 10',
  'redhat linux : 16',
  'linux : 14',
  'vmware : 7',
  'framework : 1',
  'administration : 6',
  'project : 13',
  'migration : 13',
  'tools : 3',
  'database : 29',
  'backup : 7',
  'recovery : 5',
  'production : 15',
  'development : 21',
  'production support : 5',
  'ibm : 6',
  'aix : 7',
  'servers : 41',
  'azure : 1',
  'oracle 12c : 9',
  'rac : 29',
  'node : 1',
  'security : 6',
  'installation : 8',
  'soa : 2',
  'performance monitoring : 2',
  'windows : 7',
  'database management : 1',
  'management : 24',
  'sales : 2',
  'oracle rac : 6',
  'symantec : 2',
  'automated : 1',
  'configuration : 16',
  'sqlserver : 5',
  'training : 5',
  'sql : 12',
  'rackspace : 2',
  'aws : 5',
  'oracle database : 6',
  'warehouse : 2',
  'aws cloud : 1',
  'cloud computing : 2',
  'informatica : 4',
  'etl : 3',
  'process : 1',
  'grid : 9',
  'toad : 6',
  'sql tuning : 5',
  'windows servers : 4',
  'red hat : 2',
  'rest : 1',
  'unix scripting : 4',
  'scripting : 6',
  'powershell : 1',
  'elastic search : 2',
  'git : 2',
  'jira : 2',
  'confluence : 2',
  'logstash : 2',
  'global : 3',
  'ovm : 1',
  'ux : 1',
  'oracle sql : 2',
  'sql server : 1',
  'high availability : 4',
  'data quality : 2',
  'audit : 4',
  'regulations : 1',
  'data architect : 1',
  'architect : 1',
  'netbackup : 3',
  'analysis : 11',
  'sql dba : 1',
  'solaris : 10',
  'cloud : 1',
  'weblogic : 2',
  'apex : 2',
  'performance tuning : 2',
  'data migration : 1',
  'tde : 1',
  'data encryption : 1',
  'encryption : 1',
  'perl : 2',
  'shell : 2',
  'virtualization : 3',
  'netapp : 1',
  'asm : 9',
  'stored procedures : 1',
  'tables : 1',
  'pci : 1',
  'security clearance : 5',
  'oracle fusion : 2',
  'fusion : 1',
  'scheduling : 1',
  'configuration management : 1',
  'remedy : 2',
  'incident management : 1',
  'sla : 2',
  'install : 1',
  'p : 1',
  'hardware : 1',
  'engineering : 1',
  'ms sql : 1',
  'bi : 2',
  'dw : 2',
  'siebel : 4',
  'crm : 1',
  'documentation : 2',
  'unix shell : 1',
  'shell scripting : 1',
  'mis : 1',
  'data capture : 2',
  'similar : 1',
  'pricing : 1',
  'forecasting : 1',
  'capital market : 1',
  'market research : 1',
  'research : 1',
  'trading : 1',
  'operations : 2',
  'database tuning : 2',
  'trouble shooting : 2',
  'statistics : 1',
  'queries : 1',
  'reports : 1',
  'cerner : 1',
  'emc : 3',
  'data warehousing : 1',
  'disaster recovery : 3',
  'reporting : 2',
  'veritas : 4',
  'windows server : 2',
  'credit risk : 1',
  'risk assessment : 1',
  'general ledger : 2',
  'peoplesoft : 1',
  'sybase : 1',
  'software : 9',
  'enterprise manager : 2',
  'pl : 2',
  'system administration : 3',
  'apache : 2',
  'business objects : 2',
  'business continuity : 1',
  'cluster : 6',
  'quality assurance : 1',
  'os : 1',
  'tivoli : 2',
  'network management : 7',
  'tcp : 3',
  'ip : 3',
  'cisco routers : 2',
  'routers : 2',
  'network security : 1',
  'security policies : 1',
  'policies : 1',
  'vms : 4',
  'unix administration : 1',
  'c : 2',
  'communication : 1',
  'system analysis : 2',
  'integration : 2',
  'customer support : 2',
  'gateways : 1',
  'osi : 2',
  'capacity planning : 1',
  'planning : 1',
  'protocols : 1',
  'software development : 2',
  'automation : 2',
  'systems analyst : 1',
  'technical support : 1',
  'costing : 1',
  'systems engineer : 1',
  'engineer : 1',
  'object oriented : 1',
  'design : 1',
  'embedded software : 2',
  'red hat linux : 2',
  'it : 2',
  'unix shell scripting : 1',
  'front end : 1',
  'pl/sql : 1'
)

for ($i = 0; $i -lt $keywords.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 5).Value = $keywords[$i]
}

